$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: clone per-row cell formatting (banding/borders/number formats) from matching donor rows ---
# (applied before writing values so the new cells inherit the correct alternating table-row style;
#  PasteSpecial always materialises a blank cell at every column of the destination rectangle, even
#  where the donor row had no cell at all (e.g. the unused M or N quiz-answer column), so we clear
#  that one phantom cell back out immediately after each paste.)
$ws.Range("A514:N514").Copy()
$ws.Range("A518:N518").PasteSpecial(-4122)
$ws.Range("M518").Clear()
$ws.Range("A515:N515").Copy()
$ws.Range("A519:N519").PasteSpecial(-4122)
$ws.Range("N519").Clear()
$ws.Range("A516:N516").Copy()
$ws.Range("A520:N520").PasteSpecial(-4122)
$ws.Range("N520").Clear()
$ws.Range("A515:N515").Copy()
$ws.Range("A521:N521").PasteSpecial(-4122)
$ws.Range("N521").Clear()
$ws.Range("A516:N516").Copy()
$ws.Range("A522:N522").PasteSpecial(-4122)
$ws.Range("N522").Clear()
$ws.Range("A513:N513").Copy()
$ws.Range("A523:N523").PasteSpecial(-4122)
$ws.Range("M523").Clear()
$ws.Range("A516:N516").Copy()
$ws.Range("A524:N524").PasteSpecial(-4122)
$ws.Range("N524").Clear()
$ws.Range("A513:N513").Copy()
$ws.Range("A525:N525").PasteSpecial(-4122)
$ws.Range("M525").Clear()
$ws.Range("A516:N516").Copy()
$ws.Range("A526:N526").PasteSpecial(-4122)
$ws.Range("N526").Clear()
$ws.Range("A515:N515").Copy()
$ws.Range("A527:N527").PasteSpecial(-4122)
$ws.Range("N527").Clear()
$ws.Range("A516:N516").Copy()
$ws.Range("A528:N528").PasteSpecial(-4122)
$ws.Range("N528").Clear()
$ws.Range("A513:N513").Copy()
$ws.Range("A529:N529").PasteSpecial(-4122)
$ws.Range("M529").Clear()
$ws.Range("A514:N514").Copy()
$ws.Range("A530:N530").PasteSpecial(-4122)
$ws.Range("M530").Clear()
$ws.Range("A513:N513").Copy()
$ws.Range("A531:N531").PasteSpecial(-4122)
$ws.Range("M531").Clear()
$ws.Range("A514:N514").Copy()
$ws.Range("A532:N532").PasteSpecial(-4122)
$ws.Range("M532").Clear()
$ws.Range("A515:N515").Copy()
$ws.Range("A533:N533").PasteSpecial(-4122)
$ws.Range("N533").Clear()
$ws.Range("A516:N516").Copy()
$ws.Range("A534:N534").PasteSpecial(-4122)
$ws.Range("N534").Clear()
$ws.Range("A517:N517").Copy()
$ws.Range("A535:N535").PasteSpecial(-4122)
[void]$excel.CutCopyMode
$excel.CutCopyMode = 0

# --- Step 2: write the new Google-Form response rows (518-535) ---
$ws.Range("A518").Value = 45569.922595115742
$ws.Range("B518").Value = "goeunsue@naver.com"
$ws.Range("C518").Value = "경영대학"
$ws.Range("D518").Value = 20242907
$ws.Range("E518").Value = "고은수"
$ws.Range("F518").Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Range("G518").Value = 0.1
$ws.Range("H518").Value = "6:4"
$ws.Range("I518").Value = "20분의 1"
$ws.Range("J518").Value = "20만호, 69만명"
$ws.Range("K518").Value = "충청"
$ws.Range("L518").Value = "Black"
$ws.Range("N518").Value = "모름/무응답"

$ws.Range("A519").Value = 45569.951329942131
$ws.Range("B519").Value = "taewon16@naver.com"
$ws.Range("C519").Value = "빅데이터학과"
$ws.Range("D519").Value = 20195158
$ws.Range("E519").Value = "류태원"
$ws.Range("F519").Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Range("G519").Value = 0.1
$ws.Range("H519").Value = "6:4"
$ws.Range("I519").Value = "15분의 1"
$ws.Range("J519").Value = "20만호, 69만명"
$ws.Range("K519").Value = "충청"
$ws.Range("L519").Value = "Red"
$ws.Range("M519").Value = "반대한다."

$ws.Range("A520").Value = 45569.960544363421
$ws.Range("B520").Value = "jb9517asd@naver.com"
$ws.Range("C520").Value = "소프트웨어학부"
$ws.Range("D520").Value = 20245109
$ws.Range("E520").Value = "곽우주"
$ws.Range("F520").Value = "과전법 체제에서 전국 토지를 세 등급으로 나누고 실제 수확량을 확인하여 징수하였다."
$ws.Range("G520").Value = 0.7
$ws.Range("H520").Value = "7:3"
$ws.Range("I520").Value = "10분의 1"
$ws.Range("J520").Value = "130만호, 5백만명"
$ws.Range("K520").Value = "경기"
$ws.Range("L520").Value = "Red"
$ws.Range("M520").Value = "반대한다."

$ws.Range("A521").Value = 45569.967517025463
$ws.Range("B521").Value = "1202kge@naver.com"
$ws.Range("C521").Value = "사회학과"
$ws.Range("D521").Value = 20242205
$ws.Range("E521").Value = "김가은"
$ws.Range("F521").Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Range("G521").Value = 0.7
$ws.Range("H521").Value = "4:6"
$ws.Range("I521").Value = "10분의 1"
$ws.Range("J521").Value = "44만호, 153만명"
$ws.Range("K521").Value = "전라"
$ws.Range("L521").Value = "Red"
$ws.Range("M521").Value = "근로시간과 휴무를 유연하게 조정할 수 있어 찬성한다."

$ws.Range("A522").Value = 45569.973177731481
$ws.Range("B522").Value = "sowon051125@naver.com"
$ws.Range("C522").Value = "데이터사이언스학부"
$ws.Range("D522").Value = 20243238
$ws.Range("E522").Value = "이소원"
$ws.Range("F522").Value = "‘조(租)’는 공전(公田)의 경작자가 국고에 상납하는 지대 또는 사전(私田)의 경작자가 전주에게 바치는 지대를 뜻한다."
$ws.Range("G522").Value = 0.3
$ws.Range("H522").Value = "5:5"
$ws.Range("I522").Value = "20분의 1"
$ws.Range("J522").Value = "15만호,  32만명"
$ws.Range("K522").Value = "경상"
$ws.Range("L522").Value = "Red"
$ws.Range("M522").Value = "근로시간과 휴무를 유연하게 조정할 수 있어 찬성한다."

$ws.Range("A523").Value = 45569.982878310184
$ws.Range("B523").Value = "dlxotjq27@gmail.com"
$ws.Range("C523").Value = "경영학과"
$ws.Range("D523").Value = 20213023
$ws.Range("E523").Value = "이태섭"
$ws.Range("F523").Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Range("G523").Value = 0.1
$ws.Range("H523").Value = "6:4"
$ws.Range("I523").Value = "20분의 1"
$ws.Range("J523").Value = "20만호, 69만명"
$ws.Range("K523").Value = "충청"
$ws.Range("L523").Value = "Black"
$ws.Range("N523").Value = "노동자가 과도한 연장근로를 받을 수 있어 반대한다."

$ws.Range("A524").Value = 45569.984793599535
$ws.Range("B524").Value = "jytoto33@naver.com"
$ws.Range("C524").Value = "언어청각학부"
$ws.Range("D524").Value = 20243912
$ws.Range("E524").Value = "김지윤"
$ws.Range("F524").Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Range("G524").Value = 0.3
$ws.Range("H524").Value = "3:7"
$ws.Range("I524").Value = "30분의 1"
$ws.Range("J524").Value = "44만호, 153만명"
$ws.Range("K524").Value = "평안"
$ws.Range("L524").Value = "Red"
$ws.Range("M524").Value = "반대한다."

$ws.Range("A525").Value = 45570.011931076384
$ws.Range("B525").Value = "andy041001@naver.com"
$ws.Range("C525").Value = "러시아학과"
$ws.Range("D525").Value = 20231720
$ws.Range("E525").Value = "이형범"
$ws.Range("F525").Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Range("G525").Value = 0.9
$ws.Range("H525").Value = "4:6"
$ws.Range("I525").Value = "10분의 1"
$ws.Range("J525").Value = "44만호, 153만명"
$ws.Range("K525").Value = "전라"
$ws.Range("L525").Value = "Black"
$ws.Range("N525").Value = "노동자가 과도한 연장근로를 받을 수 있어 반대한다."

$ws.Range("A526").Value = 45570.025935844911
$ws.Range("B526").Value = "hyerim0v0@gmail.com"
$ws.Range("C526").Value = "일본학과"
$ws.Range("D526").Value = 20231630
$ws.Range("E526").Value = "전혜림"
$ws.Range("F526").Value = "실제로 현장에 나가서 수확량을 파악하고 등급을 매기는 답험(踏驗)을 하였다."
$ws.Range("G526").Value = 0.3
$ws.Range("H526").Value = "6:4"
$ws.Range("I526").Value = "15분의 1"
$ws.Range("J526").Value = "20만호, 69만명"
$ws.Range("K526").Value = "경상"
$ws.Range("L526").Value = "Red"
$ws.Range("M526").Value = "반대한다."

$ws.Range("A527").Value = 45570.040854652776
$ws.Range("B527").Value = "kby5432@naver.com"
$ws.Range("C527").Value = "법학과"
$ws.Range("D527").Value = 20192737
$ws.Range("E527").Value = "윤경빈"
$ws.Range("F527").Value = "‘조(租)’는 공전(公田)의 경작자가 국고에 상납하는 지대 또는 사전(私田)의 경작자가 전주에게 바치는 지대를 뜻한다."
$ws.Range("G527").Value = 0.1
$ws.Range("H527").Value = "7:3"
$ws.Range("I527").Value = "15분의 1"
$ws.Range("J527").Value = "44만호, 153만명"
$ws.Range("K527").Value = "평안"
$ws.Range("L527").Value = "Red"
$ws.Range("M527").Value = "반대한다."

$ws.Range("A528").Value = 45570.045399247683
$ws.Range("B528").Value = "jign1106@naver.com"
$ws.Range("C528").Value = "간호학과"
$ws.Range("D528").Value = 20246289
$ws.Range("E528").Value = "지은총"
$ws.Range("F528").Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Range("G528").Value = 0.1
$ws.Range("H528").Value = "6:4"
$ws.Range("I528").Value = "20분의 1"
$ws.Range("J528").Value = "20만호, 69만명"
$ws.Range("K528").Value = "충청"
$ws.Range("L528").Value = "Red"
$ws.Range("M528").Value = "반대한다."

$ws.Range("A529").Value = 45570.080389097224
$ws.Range("B529").Value = "kt433@naver.com"
$ws.Range("C529").Value = "사회복지학과"
$ws.Range("D529").Value = 20222361
$ws.Range("E529").Value = "주혜린"
$ws.Range("F529").Value = "‘세(稅)’는 사전의 소유자가 국가에 상납하는 지대를 뜻한다."
$ws.Range("G529").Value = 0.1
$ws.Range("H529").Value = "7:3"
$ws.Range("I529").Value = "10분의 1"
$ws.Range("J529").Value = "20만호, 69만명"
$ws.Range("K529").Value = "충청"
$ws.Range("L529").Value = "Black"
$ws.Range("N529").Value = "노동자가 과도한 연장근로를 받을 수 있어 반대한다."

$ws.Range("A530").Value = 45570.130003969913
$ws.Range("B530").Value = "misunhong0707@gmail.com"
$ws.Range("C530").Value = "융합과학수사학과"
$ws.Range("D530").Value = 20246940
$ws.Range("E530").Value = "홍미선"
$ws.Range("F530").Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Range("G530").Value = 0.1
$ws.Range("H530").Value = "6:4"
$ws.Range("I530").Value = "20분의 1"
$ws.Range("J530").Value = "20만호, 69만명"
$ws.Range("K530").Value = "경상"
$ws.Range("L530").Value = "Black"
$ws.Range("N530").Value = "노동자가 과도한 연장근로를 받을 수 있어 반대한다."

$ws.Range("A531").Value = 45570.133293969906
$ws.Range("B531").Value = "ziva0726@naver.com"
$ws.Range("C531").Value = "심리학과"
$ws.Range("D531").Value = 20212104
$ws.Range("E531").Value = "김소현"
$ws.Range("F531").Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Range("G531").Value = 0.3
$ws.Range("H531").Value = "6:4"
$ws.Range("I531").Value = "30분의 1"
$ws.Range("J531").Value = "20만호, 69만명"
$ws.Range("K531").Value = "전라"
$ws.Range("L531").Value = "Black"
$ws.Range("N531").Value = "노동자가 과도한 연장근로를 받을 수 있어 반대한다."

$ws.Range("A532").Value = 45570.147500289357
$ws.Range("B532").Value = "sujdiamond@gmail.com"
$ws.Range("C532").Value = "바이오메디컬학과"
$ws.Range("D532").Value = 20243627
$ws.Range("E532").Value = "심유진"
$ws.Range("F532").Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Range("G532").Value = 0.1
$ws.Range("H532").Value = "6:4"
$ws.Range("I532").Value = "20분의 1"
$ws.Range("J532").Value = "20만호, 69만명"
$ws.Range("K532").Value = "충청"
$ws.Range("L532").Value = "Black"
$ws.Range("N532").Value = "찬성한다."

$ws.Range("A533").Value = 45570.168118692134
$ws.Range("B533").Value = "20182346@hallym.ac.kr"
$ws.Range("C533").Value = "사회복지학부"
$ws.Range("D533").Value = 20182346
$ws.Range("E533").Value = "이용재"
$ws.Range("F533").Value = "‘세(稅)’는 사전의 소유자가 국가에 상납하는 지대를 뜻한다."
$ws.Range("G533").Value = 0.1
$ws.Range("H533").Value = "7:3"
$ws.Range("I533").Value = "20분의 1"
$ws.Range("J533").Value = "130만호, 5백만명"
$ws.Range("K533").Value = "전라"
$ws.Range("L533").Value = "Red"
$ws.Range("M533").Value = "반대한다."

$ws.Range("A534").Value = 45570.174523958332
$ws.Range("B534").Value = "dncks5343@naver.com"
$ws.Range("C534").Value = "언어청각학부"
$ws.Range("D534").Value = 20243973
$ws.Range("E534").Value = "장우찬"
$ws.Range("F534").Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Range("G534").Value = 0.1
$ws.Range("H534").Value = "6:4"
$ws.Range("I534").Value = "20분의 1"
$ws.Range("J534").Value = "20만호, 69만명"
$ws.Range("K534").Value = "충청"
$ws.Range("L534").Value = "Red"
$ws.Range("M534").Value = "근로시간과 휴무를 유연하게 조정할 수 있어 찬성한다."

$ws.Range("A535").Value = 45570.24906371528
$ws.Range("B535").Value = "leedongyoung797@gmail.com"
$ws.Range("C535").Value = "언어청각학부"
$ws.Range("D535").Value = 20243934
$ws.Range("E535").Value = "이동영"
$ws.Range("F535").Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Range("G535").Value = 0.3
$ws.Range("H535").Value = "6:4"
$ws.Range("I535").Value = "10분의 1"
$ws.Range("J535").Value = "20만호, 69만명"
$ws.Range("K535").Value = "전라"
$ws.Range("L535").Value = "Red"
$ws.Range("M535").Value = "근로시간과 휴무를 유연하게 조정할 수 있어 찬성한다."

# --- Step 3: row 517 is no longer the last table row, so drop its old empty placeholder cell N517 ---
$ws.Range("N517").Clear()

# --- Step 4: expand the table (ListObject) to the new extent ---
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:N535"))

# --- Step 5: move the saved selection to match the author's last cursor position ---
[void]$ws.Range("C542").Select()
